# Update the "Pais" (countries) dashboard with the latest COVID-19 figures
# and refresh the "last updated" timestamp.
#
# Columns: A=Pais, B=Casos totales, C=Nuevos casos, D=Casos activos,
#          E=Recuperados, F=Casos criticos, G=Muertes hoy, H=Muertes

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp banner (A1)
$ws.Range("A1").Value = "Datos actualizados a 25 de Octubre de 2020 a las 05:09"

# India (row 5)
$ws.Range("B5").Value = 7864811
$ws.Range("C5").Value = 919
$ws.Range("D5").Value = 7078123
$ws.Range("E5").Value = 668121

# Alemania (row 20)
$ws.Range("B20").Value = 429181
$ws.Range("C20").Value = 1373
$ws.Range("D20").Value = 317000
$ws.Range("E20").Value = 102070

# Belgica (row 29)
$ws.Range("B29").Value = 305409
$ws.Range("C29").Value = 17709
$ws.Range("D29").Value = 22878
$ws.Range("E29").Value = 271794
$ws.Range("G29").Value = 79
$ws.Range("H29").Value = 10737

# Honduras overtakes Bielorrusia, so the two rows swap places (54/55)
# Row 54: now Honduras (was Bielorrusia)
$ws.Range("A54").Value = "Honduras"
$ws.Range("B54").Value = 92724
$ws.Range("C54").Value = 842
$ws.Range("D54").Value = 37866
$ws.Range("E54").Value = 52241
$ws.Range("G54").Value = 5
$ws.Range("H54").Value = 2617

# Row 55: now Bielorrusia (was Honduras)
$ws.Range("A55").Value = "Bielorrusia"
$ws.Range("B55").Value = 91978
$ws.Range("D55").Value = 82670
$ws.Range("E55").Value = 8355
$ws.Range("H55").Value = 953

# Sudan (row 103)
$ws.Range("B103").Value = 13742
$ws.Range("E103").Value = 6141
$ws.Range("H103").Value = 837

# Tailandia (row 149)
$ws.Range("B149").Value = 3736
$ws.Range("C149").Value = 5
$ws.Range("D149").Value = 3530
$ws.Range("E149").Value = 147

# Belice (row 153)
$ws.Range("B153").Value = 3106
$ws.Range("C153").Value = 56
$ws.Range("D153").Value = 1921
$ws.Range("E153").Value = 1137
$ws.Range("G153").Value = 2
$ws.Range("H153").Value = 48

# Burkina Faso (row 157)
$ws.Range("B157").Value = 2444
$ws.Range("D157").Value = 1997
$ws.Range("E157").Value = 382

# Butan overtakes Mongolia, so the two rows swap places (187/188)
# Row 187: now Butan (was Mongolia)
$ws.Range("A187").Value = "Butan"
$ws.Range("B187").Value = 340
$ws.Range("C187").Value = 4
$ws.Range("D187").Value = 306
$ws.Range("E187").Value = 34

# Row 188: now Mongolia (was Butan)
$ws.Range("A188").Value = "Mongolia"
$ws.Range("B188").Value = 338
$ws.Range("C188").Value = 1
$ws.Range("D188").Value = 312
$ws.Range("E188").Value = 26

# Montserrat and Islas Malvinas swap places (216/217)
# Row 216: now Montserrat (was Islas Malvinas)
$ws.Range("A216").Value = "Montserrat"
$ws.Range("B216").Value = 12
$ws.Range("D216").Value = 12
$ws.Range("H216").Value = 1

# Row 217: now Islas Malvinas (was Montserrat)
$ws.Range("A217").Value = "Islas Malvinas"
$ws.Range("D217").Value = 13
$ws.Range("H217").Value = 0
